# Update column G ("K" / Strikeouts) values for rows 2-16 per regenerated save_data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 2
    3  = 0
    4  = 0
    5  = 0
    6  = 2
    7  = 3
    8  = 2
    9  = 0
    10 = 2
    11 = 0
    12 = 1
    13 = 0
    14 = 1
    15 = 1
    16 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
